$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target range to Text format so numeric-looking strings
# (e.g. "602.19") are kept verbatim instead of being parsed into floats,
# then restore the original (default) style so no stray formatting is left behind.
$dataRange = $ws.Range("D2:E51")
$originalStyle = $ws.Cells.Item(2, 4).Style
$dataRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '65.807.77'
$ws.Cells.Item(2, 5).Value = '  +0.73%  '
$ws.Cells.Item(3, 4).Value = '2.680.44'
$ws.Cells.Item(3, 5).Value = '  +0.93%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).Value = '602.19'
$ws.Cells.Item(5, 5).Value = '  -0.66%  '
$ws.Cells.Item(6, 4).Value = '156.46'
$ws.Cells.Item(6, 5).Value = '  -0.83%  '
$ws.Cells.Item(7, 5).Value = '  +0.00%  '
$ws.Cells.Item(8, 4).Value = '0.588'
$ws.Cells.Item(8, 5).Value = '  -0.15%  '
$ws.Cells.Item(9, 5).Value = '  -0.12%  '
$ws.Cells.Item(10, 4).Value = '5.93'
$ws.Cells.Item(10, 5).Value = '  +1.71%  '
$ws.Cells.Item(11, 5).Value = '  -3.57%  '
$ws.Cells.Item(12, 5).Value = '  +0.14%  '
$ws.Cells.Item(13, 4).Value = '29.53'
$ws.Cells.Item(13, 5).Value = '  -0.17%  '
$ws.Cells.Item(14, 5).Value = '  +6.51%  '
$ws.Cells.Item(15, 4).Value = '3.159.02'
$ws.Cells.Item(15, 5).Value = '  +0.81%  '
$ws.Cells.Item(16, 4).Value = '65.613.14'
$ws.Cells.Item(16, 5).Value = '  +0.74%  '
$ws.Cells.Item(17, 4).Value = '2.689.88'
$ws.Cells.Item(17, 5).Value = '  +1.30%  '
$ws.Cells.Item(18, 4).Value = '12.62'
$ws.Cells.Item(18, 5).Value = '  -1.31%  '
$ws.Cells.Item(19, 4).Value = '4.83'
$ws.Cells.Item(19, 5).Value = '  -1.88%  '
$ws.Cells.Item(20, 4).Value = '7.61'
$ws.Cells.Item(20, 5).Value = '  +3.21%  '
$ws.Cells.Item(21, 4).Value = '352.28'
$ws.Cells.Item(21, 5).Value = '  -2.23%  '
$ws.Cells.Item(23, 4).Value = '70.09'
$ws.Cells.Item(23, 5).Value = '  +1.13%  '
$ws.Cells.Item(24, 5).Value = '  +7.27%  '
$ws.Cells.Item(25, 4).Value = '9.82'
$ws.Cells.Item(25, 5).Value = '  +2.34%  '
$ws.Cells.Item(26, 5).Value = '  -5.25%  '
$ws.Cells.Item(27, 5).Value = '  -1.87%  '
$ws.Cells.Item(28, 5).Value = '  +2.32%  '
$ws.Cells.Item(29, 4).Value = '8.19'
$ws.Cells.Item(29, 5).Value = '  -0.86%  '
$ws.Cells.Item(30, 5).Value = '  -0.09%  '
$ws.Cells.Item(31, 5).Value = '  -2.73%  '
$ws.Cells.Item(32, 4).Value = '531.02'
$ws.Cells.Item(32, 5).Value = '  -4.23%  '
$ws.Cells.Item(33, 5).Value = '  -3.82%  '
$ws.Cells.Item(34, 4).Value = '6.56'
$ws.Cells.Item(34, 5).Value = '  +2.12%  '
$ws.Cells.Item(35, 5).Value = '  -3.76%  '
$ws.Cells.Item(36, 4).Value = '0.427'
$ws.Cells.Item(36, 5).Value = '  -1.66%  '
$ws.Cells.Item(37, 4).Value = '20.51'
$ws.Cells.Item(37, 5).Value = '  -0.55%  '
$ws.Cells.Item(38, 4).Value = '160.96'
$ws.Cells.Item(38, 5).Value = '  -1.62%  '
$ws.Cells.Item(39, 5).Value = '  +0.02%  '
$ws.Cells.Item(40, 4).Value = '1.97'
$ws.Cells.Item(40, 5).Value = '  -2.17%  '
$ws.Cells.Item(41, 5).Value = '  -0.01%  '
$ws.Cells.Item(42, 4).Value = '42.27'
$ws.Cells.Item(42, 5).Value = '  -0.57%  '
$ws.Cells.Item(43, 4).Value = '166.16'
$ws.Cells.Item(43, 5).Value = '  -0.99%  '
$ws.Cells.Item(44, 4).Value = '4.11'
$ws.Cells.Item(44, 5).Value = '  -2.42%  '
$ws.Cells.Item(45, 4).Value = '0.0621'
$ws.Cells.Item(45, 5).Value = '  +0.00%  '
$ws.Cells.Item(46, 4).Value = '23.16'
$ws.Cells.Item(46, 5).Value = '  +0.31%  '
$ws.Cells.Item(47, 5).Value = '  -3.86%  '
$ws.Cells.Item(48, 5).Value = '  -0.56%  '
$ws.Cells.Item(49, 5).Value = '  -1.01%  '
$ws.Cells.Item(50, 4).Value = '20.33'
$ws.Cells.Item(50, 5).Value = '  +2.87%  '
$ws.Cells.Item(51, 4).Value = '0.0988'
$ws.Cells.Item(51, 5).Value = '  +0.21%  '

$dataRange.Style = $originalStyle

